$d = $word.ActiveDocument

# Find the paragraph holding the site footer/copyright notice
# ("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
# pages. Original theme under Creative Commons Attribution"). In this
# document it is preceded by two blank "spacer" paragraphs: an empty
# Normal paragraph and an empty paragraph carrying a page break. The
# whole block (the two spacers plus the footer paragraph) is removed,
# while the requisite line above it and the trailing blank/page-break
# paragraphs that originally followed it are left untouched.

$footerIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Powered by Jekyll*") {
        $footerIndex = $i
        break
    }
}

if ($footerIndex -ge 3) {
    $firstPara = $d.Paragraphs.Item($footerIndex - 2)
    $lastPara  = $d.Paragraphs.Item($footerIndex)
    $rangeToDelete = $d.Range($firstPara.Range.Start, $lastPara.Range.End)
    $rangeToDelete.Delete()
}
